# Update "Förändrad" (Changed) date column C for rows 2-11.
# Each cell's underlying serial date value moves from 45233 (2023-11-03)
# to 45243 (2023-11-13), i.e. +10 days, while keeping the existing
# yyyy-mm-dd style/format untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
